$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.048.34"
$ws.Range("E2").Value = "  -2.96%  "

$ws.Range("D3").Value = "2.577.44"
$ws.Range("E3").Value = "  -2.21%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'534.35"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").Value = "'141.28"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  +3.73%  "

$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("D10").Value = "'0.0995"
$ws.Range("E10").Value = "  -3.27%  "

$ws.Range("E11").Value = "  +2.72%  "

$ws.Range("E12").Value = "  -2.06%  "

$ws.Range("D13").Value = "3.032.61"
$ws.Range("E13").Value = "  -2.46%  "

$ws.Range("D14").Value = "58.000.01"
$ws.Range("E14").Value = "  -2.91%  "

$ws.Range("E15").Value = "  -1.34%  "

$ws.Range("D16").Value = "2.559.67"
$ws.Range("E16").Value = "  -3.38%  "

$ws.Range("E17").Value = "  -2.70%  "

$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "'333.69"
$ws.Range("E19").Value = "  -2.75%  "

$ws.Range("D20").Value = "'10.04"
$ws.Range("E20").Value = "  -1.67%  "

$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  -3.97%  "

$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'66.69"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "'0.159"
$ws.Range("E26").Value = "  -4.19%  "

$ws.Range("D27").Value = "'7.02"
$ws.Range("E27").Value = "  -3.52%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").Value = "0.0₃0724"
$ws.Range("E29").Value = "  -3.72%  "

$ws.Range("E30").Value = "  -2.32%  "

$ws.Range("D31").Value = "'155.69"
$ws.Range("E31").Value = "  +3.09%  "

$ws.Range("D32").Value = "'5.86"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").Value = "'18.82"
$ws.Range("E33").Value = "  -0.51%  "

$ws.Range("E34").Value = "  -3.41%  "

$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("E36").Value = "  -4.01%  "

$ws.Range("D37").Value = "'0.831"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("E38").Value = "  -2.52%  "

$ws.Range("E39").Value = "  -3.90%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "'281.01"
$ws.Range("E41").Value = "  -3.68%  "

$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0951"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'10.64"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("E47").Value = "  -2.30%  "

$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").Value = "1.906.15"
$ws.Range("E49").Value = "  -3.20%  "

$ws.Range("D50").Value = "'17.76"
$ws.Range("E50").Value = "  -4.25%  "

$ws.Range("E51").Value = "  -3.41%  "
